# Generate Report for Handback
# Applies the "handback" pass over the localization-status workbook:
#  - Overview status cells flip from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Each language sheet (zh-cn / de-de) gets its "Latest Target File" (hyperlink),
#    "Latest Handback File" and "Latest Handback DateTime" columns filled in for both rows
#  - A few columns are widened to fit the newly-populated content

$wb = $excel.ActiveWorkbook

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3be06001de85acdc2d77eae2ace849a9b375c98a/e2e/113e5753-58ac-480c-a23a-eb8aa64611be.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3be06001de85acdc2d77eae2ace849a9b375c98a/e2e/ced3542b-72a3-4ee4-8eb6-502ddcdb43db.md"

$name1 = "113e5753-58ac-480c-a23a-eb8aa64611be.md"
$name2 = "ced3542b-72a3-4ee4-8eb6-502ddcdb43db.md"

# ---------------------------------------------------------------------------
# Overview sheet: both locale status columns report the handback is in sync
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $urlMd1, [Type]::Missing, [Type]::Missing, $name1)
$zhcn.Range("J2").Value = "113e5753-58ac-480c-a23a-eb8aa64611be.21fbd084e97c60b4c509ef8891fd9b038157b1f6.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-15 08:45:34"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $urlMd2, [Type]::Missing, [Type]::Missing, $name2)
$zhcn.Range("J3").Value = "ced3542b-72a3-4ee4-8eb6-502ddcdb43db.7a007b8ce27983848334a52af127f922646fb17d.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-15 08:45:34"

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Hyperlinks.Add($dede.Range("I2"), $urlMd1, [Type]::Missing, [Type]::Missing, $name1)
$dede.Range("J2").Value = "113e5753-58ac-480c-a23a-eb8aa64611be.21fbd084e97c60b4c509ef8891fd9b038157b1f6.de-de.xlf"
$dede.Range("K2").Value = "2016-08-15 08:45:41"

$dede.Hyperlinks.Add($dede.Range("I3"), $urlMd2, [Type]::Missing, [Type]::Missing, $name2)
$dede.Range("J3").Value = "ced3542b-72a3-4ee4-8eb6-502ddcdb43db.7a007b8ce27983848334a52af127f922646fb17d.de-de.xlf"
$dede.Range("K3").Value = "2016-08-15 08:45:41"

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
